# Add new columns I (I0) and J (IF) with per-row values, mirroring the
# style used for existing data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
# Match header styling used by the other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-38
# Each entry: row, I value, J value
$values = @(
    ,@(2,  9, 9)
    ,@(3,  6, 6)
    ,@(4,  6, 6)
    ,@(5,  6, 6)
    ,@(6,  7, 7)
    ,@(7,  6, 6)
    ,@(8,  7, 7)
    ,@(9,  7, 7)
    ,@(10, 5, 6)
    ,@(11, 8, 8)
    ,@(12, 7, 8)
    ,@(13, 7, 8)
    ,@(14, 7, 8)
    ,@(15, 6, 7)
    ,@(16, 4, 4)
    ,@(17, 5, 5)
    ,@(18, 5, 6)
    ,@(19, 5, 5)
    ,@(20, 8, 8)
    ,@(21, 8, 8)
    ,@(22, 7, 8)
    ,@(23, 6, 7)
    ,@(24, 8, 8)
    ,@(25, 1, 1)
    ,@(26, 5, 5)
    ,@(27, 7, 8)
    ,@(28, 8, 8)
    ,@(29, 5, 5)
    ,@(30, 5, 5)
    ,@(31, 5, 5)
    ,@(32, 6, 6)
    ,@(33, 4, 4)
    ,@(34, 8, 8)
    ,@(35, 8, 8)
    ,@(36, 9, 9)
    ,@(37, 8, 8)
    ,@(38, 5, 5)
)

foreach ($entry in $values) {
    $row = $entry[0]
    $ws.Cells.Item($row, 9).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
}
